# Matrice di tracciabilità.xlsx - apply author's edit
#
# Summary of the change (derived from the OOXML diff):
#  - Status column (E7:E16) values change from "N" (non implementato) to
#    "Im" (implementato) -> the shared string "N" becomes unused/removed
#    and a new shared string "Im" is introduced.
#  - The "TC" column had two sub-columns (M = SDD/TCS/TC, N = TESTING/TC).
#    For rows 7, 8 and 15 the N-column placeholder "\" is replaced by the
#    same Test-Case id already present in the M column on that row
#    (TC_1.1, TC_2.2, TC_16.16 respectively). N7 additionally loses its
#    special underlined-font style, now matching the plain style used by
#    the rest of the data cells (same look as N8/N15/M7).
#  - Minor cosmetic view changes: zoom level and current selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Status column: "N" -> "Im" for the requirement rows ------------------
$ws.Range("E7:E16").Value = "Im"

# --- TESTING/TC column (N): fill in the same TC id as the SDD/TCS column (M) ---
$ws.Range("N7").Value = $ws.Range("M7").Value2
$ws.Range("N8").Value = $ws.Range("M8").Value2
$ws.Range("N15").Value = $ws.Range("M15").Value2

# N7 previously carried a distinct (underlined) style; align it with the
# normal data-cell formatting used elsewhere in the table (e.g. M7/N8).
$ws.Range("M7").Copy()
$ws.Range("N7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Cosmetic: view zoom + current selection, as left by the author -------
$ws.Select()
$excel.ActiveWindow.Zoom = 130
$ws.Range("C19").Select()

# --- Cosmetic: a handful of column widths were tweaked --------------------
$ws.Columns.Item(3).ColumnWidth = 39
$ws.Columns.Item(4).ColumnWidth = 8.7109375
$ws.Columns.Item(6).ColumnWidth = 20.5703125
$ws.Columns.Item(9).ColumnWidth = 15
$ws.Columns.Item(10).ColumnWidth = 24.7109375
$ws.Columns.Item(13).ColumnWidth = 11.28515625
$ws.Columns.Item(14).ColumnWidth = 11
